# Retraining the quarterly Elnet model
# Shift all timestamps in column A (rows 2..97) forward by 4 days,
# and update the forecasted production values in column B for the
# affected rows (22..39) to the newly retrained model's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in A2:A97 forward by 4 days (keeping time-of-day fraction)
for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $cell.Value2() + 4
}

# Updated production values (MW) produced by the retrained model
$bUpdates = @{
    22 = 1
    23 = 16
    25 = 82
    26 = 136
    27 = 198
    28 = 283
    29 = 374
    30 = 477
    31 = 587
    32 = 677
    33 = 778
    34 = 883
    35 = 959
    36 = 1037
    37 = 1099
    38 = 1150
    39 = 1231
}

foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}
